$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.717.57"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "3.465.11"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "413.78"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.24"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.725"
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  +5.41%  "
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.62"
$ws.Range("E12").Value = "  +4.21%  "
$ws.Range("E13").Value = "  -2.92%  "
$ws.Range("D14").Value = "4.017.36"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.141"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.54"
$ws.Range("E16").Value = "  -3.98%  "
$ws.Range("D17").Value = "3.468.41"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.69"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("D20").Value = "62.664.54"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "462.20"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.53"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.27"
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.30"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.69"
$ws.Range("E25").Value = "  +16.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.30"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.35"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.79"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.58"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.67"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("E33").Value = "  -2.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.78"
$ws.Range("E34").Value = "  -5.51%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  +7.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0491"
$ws.Range("E37").Value = "  -2.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.08"
$ws.Range("E38").Value = "  +4.34%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "147.44"
$ws.Range("E40").Value = "  +3.47%  "
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.320"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("E44").Value = "  +5.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.38"
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("E46").Value = "  +3.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("E47").Value = "  +12.96%  "
$ws.Range("D48").Value = "0.0₃0558"
$ws.Range("E48").Value = "  +29.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.42"
$ws.Range("E49").Value = "  -1.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.17"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("E51").Value = "  -1.76%  "
